$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting rows 28-32 down to 29-33
$ws.Rows.Item(28).Insert()

# Populate new row 28 with the new data
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44769
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112022
$ws.Cells.Item(28, 7).Value = "Arveja Verde"
$ws.Cells.Item(28, 8).Value = "Perfection"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 40000
$ws.Cells.Item(28, 12).Value = 42000
$ws.Cells.Item(28, 13).Value = 41000
$ws.Cells.Item(28, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(28, 16).Value = 1640
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
